$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q1" sheet, positioned right before "总计".
#    We copy the "2021-Q4" sheet (same column layout: 基金代码/基金名称/...)
#    so the new sheet inherits the same sheetPr/pageMargins/column layout,
#    then rename it and overwrite its data.
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")
$template.Copy($totalSheet)

$new = $wb.Worksheets.Item("2021-Q4 (2)")
$new.Name = "2022-Q1"

$ws = $wb.Worksheets.Item("2022-Q1")

# Header row
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# Data rows (fund code / numeric-looking text columns keep their original
# textual formatting via a leading apostrophe, so leading/trailing zeros
# such as "002295" or "2.30" are preserved instead of being coerced into
# numbers).
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "'166301"
$ws.Range("C2").Value = "华商新趋势优选灵活配置混合"
$ws.Range("D2").Value = "'26.96"
$ws.Range("E2").Value = "'86.39"
$ws.Range("F2").Value = "'2.30"
$ws.Range("G2").Value = "'0.6201"
$ws.Range("H2").Value = 10

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "'200012"
$ws.Range("C3").Value = "长城中小盘成长混合"
$ws.Range("D3").Value = "'12.65"
$ws.Range("E3").Value = "'84.26"
$ws.Range("F3").Value = "'2.42"
$ws.Range("G3").Value = "'0.3061"
$ws.Range("H3").Value = 5

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "'206002"
$ws.Range("C4").Value = "鹏华精选成长混合"
$ws.Range("D4").Value = "'4.48"
$ws.Range("E4").Value = "'92.68"
$ws.Range("F4").Value = "'4.78"
$ws.Range("G4").Value = "'0.2141"
$ws.Range("H4").Value = 5

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "'630008"
$ws.Range("C5").Value = "华商策略精选混合"
$ws.Range("D5").Value = "'4.71"
$ws.Range("E5").Value = "'76.99"
$ws.Range("F5").Value = "'2.33"
$ws.Range("G5").Value = "'0.1097"
$ws.Range("H5").Value = 9

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "'002295"
$ws.Range("C6").Value = "广发稳安灵活配置混合A"
$ws.Range("D6").Value = "'1.85"
$ws.Range("E6").Value = "'70.60"
$ws.Range("F6").Value = "'4.48"
$ws.Range("G6").Value = "'0.0829"
$ws.Range("H6").Value = 5

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "'011765"
$ws.Range("C7").Value = "兴银高端制造混合A"
$ws.Range("D7").Value = "'1.01"
$ws.Range("E7").Value = "'93.23"
$ws.Range("F7").Value = "'3.08"
$ws.Range("G7").Value = "'0.0311"
$ws.Range("H7").Value = 5

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "'011766"
$ws.Range("C8").Value = "兴银高端制造混合C"
$ws.Range("D8").Value = "'0.39"
$ws.Range("E8").Value = "'93.23"
$ws.Range("F8").Value = "'3.08"
$ws.Range("G8").Value = "'0.0120"
$ws.Range("H8").Value = 5

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "'005146"
$ws.Range("C9").Value = "兴银丰润灵活配置混合"
$ws.Range("D9").Value = "'0.05"
$ws.Range("E9").Value = "'93.36"
$ws.Range("F9").Value = "'3.96"
$ws.Range("G9").Value = "'0.0020"
$ws.Range("H9").Value = 3

# Row 10 does not exist in the copied template (it only had 8 data rows),
# so add it and copy the row-9 formatting (bold/centered/bordered A column
# style) down onto the new row's index cell.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "'008604"
$ws.Range("C10").Value = "广发稳安灵活配置混合C"
$ws.Range("D10").Value = "'0.02"
$ws.Range("E10").Value = "'70.60"
$ws.Range("F10").Value = "'4.48"
$ws.Range("G10").Value = "'0.0009"
$ws.Range("H10").Value = 5

# ---------------------------------------------------------------------------
# 2. Update the "总计" sheet: prepend the 2022-Q1 summary row and shift the
#    existing quarters down by one row.
# ---------------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")

# Row 6 does not exist yet - add it and clone row 5's index-cell formatting.
$tot.Range("A5").Copy()
$tot.Range("A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$tot.Range("A2").Value = 0
$tot.Range("B2").Value = "2022-Q1"
$tot.Range("C2").Value = 9
$tot.Range("D2").Value = 1.38

$tot.Range("A3").Value = 1
$tot.Range("B3").Value = "2021-Q4"
$tot.Range("C3").Value = 8
$tot.Range("D3").Value = 1.62

$tot.Range("A4").Value = 2
$tot.Range("B4").Value = "2021-Q3"
$tot.Range("C4").Value = 11
$tot.Range("D4").Value = 0.38

$tot.Range("A5").Value = 3
$tot.Range("B5").Value = "2021-Q2"
$tot.Range("C5").Value = 3
$tot.Range("D5").Value = 1.17

$tot.Range("A6").Value = 4
$tot.Range("B6").Value = "2021-Q1"
$tot.Range("C6").Value = 4
$tot.Range("D6").Value = 1.09

# Restore the workbook's originally active sheet/selection.
$wb.Worksheets.Item("2021-Q1").Activate() | Out-Null
